# Fixed typo in Nutshell diagram.
#
# 1) Update the "last edited" date field (stored literal text behind the
#    datetimeFigureOut field) on the slide master and every slide layout
#    from 3/29/2016 to 9/22/2016.
# 2) Fix the "MLib" typo on slide 1's flow chart to read "ML/MLLib".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the date placeholder text wherever it appears (slide master +
#    all custom layouts) by locating the placeholder of type "Date" and
#    rewriting its text.
# ---------------------------------------------------------------------
function Update-DatePlaceholder {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {   # ppPlaceholderDate
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "9/22/2016"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "9/22/2016"
}

# ---------------------------------------------------------------------
# 2. Fix "MLib" -> "ML/MLLib" in the Nutshell diagram on slide 1.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "MLib") {
            $tr = $shp.TextFrame.TextRange
            $tr.Text = "ML/"
            $tr.InsertAfter("MLLib") | Out-Null
        }
    }
}
